$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.08544659614563
$ws.Range("B1").Value = 2.581828832626343
$ws.Range("C1").Value = 2.684499025344849
$ws.Range("D1").Value = 3.169137954711914
$ws.Range("E1").Value = 0.7647852301597595
